$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 new rows at the top to make room for the new "Site Code" table.
$ws.Rows("1:7").Insert()

# New "Site Code" table (rows 1-5)
$ws.Range("A1").Value = "Site Code"
$ws.Range("B1").Value = "Full Name"

$ws.Range("A2").Value = "IC1"
$ws.Range("B2").Value = "Indian Creek 1"

$ws.Range("A3").Value = "IC2"
$ws.Range("B3").Value = "Indian Creek 2"

$ws.Range("A4").Value = "IC3"
$ws.Range("B4").Value = "Indian Creek 3"

$ws.Range("A5").Value = "SF"
$ws.Range("B5").Value = "Sand Flats"

# Vascular plant table header (now at row 8 after the insert)
$ws.Range("A8").Value = "Vascular Code"
$ws.Range("B8").Value = "Scientific Name"
$ws.Range("C8").Value = "Common Name"

# Biocrust table header (now at row 27 after the insert)
$ws.Range("A27").Value = "Biocrust Code"
$ws.Range("B27").Value = "Full Name"

# Update biocrust descriptions
$ws.Range("B32").Value = "Squamulose (Lichen)"
$ws.Range("B33").Value = "Collema tenax (Lichen)"

# Fit column A to the new longer header labels
$ws.Columns("A").AutoFit()

$ws.Range("B32").Select()
